$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2-22) holds one row per "origin" quarter, with a
# tapering number of forecast-horizon columns (B..H) the further down the
# sheet we go. A quarter ("2020-04-01") was missing from the original
# export because of a selection-scope bug; this script re-inserts it as
# the new row 3 and pushes every following origin-quarter row down by one,
# which is exactly what the corrected scope now produces.

# Column list used by each row before the fix (row -> last column letter).
$cols = @("A","B","C","D","E","F","G","H")

# Number of data columns (including A) present in each of the original
# rows 3..22, taken from the existing sheet (rows 3-16 have 8 columns,
# then the count shrinks by one every row down to row 22 which has 2).
$colCount = @{}
for ($r = 3; $r -le 16; $r++) { $colCount[$r] = 8 }
$colCount[17] = 7
$colCount[18] = 6
$colCount[19] = 5
$colCount[20] = 4
$colCount[21] = 3
$colCount[22] = 2

# Shift the existing rows 3..22 down to rows 4..23, working from the
# bottom up so we never overwrite a row before we've read it. Only the
# columns that actually contained data are touched, so no stray empty
# cells get created on the tapering rows.
for ($r = 22; $r -ge 3; $r--) {
    $n = $colCount[$r]
    $destRow = $r + 1
    for ($i = 0; $i -lt $n; $i++) {
        $col = $cols[$i]
        $val = $ws.Range($col + $r).Value2
        $ws.Range($col + $destRow).Value2 = $val
    }
}

# Row 23 did not exist before, so it has no style yet. Give its label
# cell (A23) the same look as every other label cell in column A (bold,
# centered, top-aligned, thin border) by copying the formatting only
# (not the value) from an existing label cell.
$ws.Range("A2").Copy()
$ws.Range("A23").PasteSpecial(-4122)

# Fill in the newly inserted row 3 with the corrected quarter label and
# its matched-error values.
$ws.Range("A3").Value2 = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value2 = 6.652313087672924
$ws.Range("C3").Value2 = -10.89834099542839
$ws.Range("D3").Value2 = -3.188944395772239
$ws.Range("E3").Value2 = -1.606031734885327
$ws.Range("F3").Value2 = -4.572318727234583
$ws.Range("G3").Value2 = -2.30828757257012
$ws.Range("H3").Value2 = -2.759441727496859
